$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused highlight/font style from C5:D8 (reverts them to default style)
$ws.Range("C5:D8").ClearFormats()

# Add the new data row for group 50433 that was added back to the dataset
$ws.Range("A9").Value = 50433
$ws.Range("B9").Value = 0.13600000000000001
$ws.Range("C9").Value = 157.732
$ws.Range("D9").Value = 30.844000000000001

# Update the active selection to match the author's saved cursor position
[void]$ws.Range("D16").Select()
